# Rename the second worksheet from "LoginTest" to "LoginFunc"
$wb = $excel.ActiveWorkbook
$loginSheet = $wb.Worksheets.Item("LoginTest")
$loginSheet.Name = "LoginFunc"

# Select cell D10 on the LoginFunc sheet (was D4)
$loginSheet.Range("D10").Select()

# Make LoginFunc the active sheet (this also sets tabSelected on it and
# clears tabSelected on the previously active sheet, and sets activeTab
# in the workbook view)
$loginSheet.Activate()
